$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-08-07 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-08-08 Thursday", 2) | Out-Null
$d.Content.Find.Execute("37+20=57", $true, $false, $false, $false, $false, $true, 1, $false, "95-60=35", 2) | Out-Null
$d.Content.Find.Execute("79-36=43", $true, $false, $false, $false, $false, $true, 1, $false, "88+5=93", 2) | Out-Null
$d.Content.Find.Execute("91-48=43", $true, $false, $false, $false, $false, $true, 1, $false, "60-3=57", 2) | Out-Null
$d.Content.Find.Execute("13+84=97", $true, $false, $false, $false, $false, $true, 1, $false, "5+58=63", 2) | Out-Null
$d.Content.Find.Execute("60+22=82", $true, $false, $false, $false, $false, $true, 1, $false, "62-6=56", 2) | Out-Null
$d.Content.Find.Execute("75-39=36", $true, $false, $false, $false, $false, $true, 1, $false, "16+20=36", 2) | Out-Null
$d.Content.Find.Execute("10+20=30", $true, $false, $false, $false, $false, $true, 1, $false, "60-28=32", 2) | Out-Null
$d.Content.Find.Execute("66+18=84", $true, $false, $false, $false, $false, $true, 1, $false, "86-2=84", 2) | Out-Null
$d.Content.Find.Execute("1+33=34", $true, $false, $false, $false, $false, $true, 1, $false, "62+18=80", 2) | Out-Null
$d.Content.Find.Execute("45-44=1", $true, $false, $false, $false, $false, $true, 1, $false, "19+15=34", 2) | Out-Null
$d.Content.Find.Execute("86-76=10", $true, $false, $false, $false, $false, $true, 1, $false, "14+20=34", 2) | Out-Null
$d.Content.Find.Execute("84-4=80", $true, $false, $false, $false, $false, $true, 1, $false, "90-51=39", 2) | Out-Null
$d.Content.Find.Execute("38-9=29", $true, $false, $false, $false, $false, $true, 1, $false, "86-42=44", 2) | Out-Null
$d.Content.Find.Execute("82-60=22", $true, $false, $false, $false, $false, $true, 1, $false, "66-21=45", 2) | Out-Null
$d.Content.Find.Execute("20-12=8", $true, $false, $false, $false, $false, $true, 1, $false, "0+12=12", 2) | Out-Null
$d.Content.Find.Execute("61-9=52", $true, $false, $false, $false, $false, $true, 1, $false, "5+37=42", 2) | Out-Null
$d.Content.Find.Execute("89-22=67", $true, $false, $false, $false, $false, $true, 1, $false, "49+35=84", 2) | Out-Null
$d.Content.Find.Execute("55-23=32", $true, $false, $false, $false, $false, $true, 1, $false, "30+9=39", 2) | Out-Null
$d.Content.Find.Execute("87-32=55", $true, $false, $false, $false, $false, $true, 1, $false, "68-28=40", 2) | Out-Null
$d.Content.Find.Execute("44+34=78", $true, $false, $false, $false, $false, $true, 1, $false, "27+35=62", 2) | Out-Null
$d.Content.Find.Execute("33+34=67", $true, $false, $false, $false, $false, $true, 1, $false, "34-3=31", 2) | Out-Null
$d.Content.Find.Execute("52+44=96", $true, $false, $false, $false, $false, $true, 1, $false, "80-24=56", 2) | Out-Null
$d.Content.Find.Execute("94-54=40", $true, $false, $false, $false, $false, $true, 1, $false, "56+0=56", 2) | Out-Null
$d.Content.Find.Execute("69-14=55", $true, $false, $false, $false, $false, $true, 1, $false, "20+61=81", 2) | Out-Null
$d.Content.Find.Execute("91-81=10", $true, $false, $false, $false, $false, $true, 1, $false, "56+12=68", 2) | Out-Null
$d.Content.Find.Execute("42+51=93", $true, $false, $false, $false, $false, $true, 1, $false, "40+37=77", 2) | Out-Null
$d.Content.Find.Execute("16+3=19", $true, $false, $false, $false, $false, $true, 1, $false, "3+20=23", 2) | Out-Null
$d.Content.Find.Execute("17+72=89", $true, $false, $false, $false, $false, $true, 1, $false, "41-6=35", 2) | Out-Null
$d.Content.Find.Execute("63+0=63", $true, $false, $false, $false, $false, $true, 1, $false, "11+27=38", 2) | Out-Null
$d.Content.Find.Execute("40-3=37", $true, $false, $false, $false, $false, $true, 1, $false, "9+4=13", 2) | Out-Null
$d.Content.Find.Execute("5+92=97", $true, $false, $false, $false, $false, $true, 1, $false, "49-36=13", 2) | Out-Null
$d.Content.Find.Execute("14-2=12", $true, $false, $false, $false, $false, $true, 1, $false, "25-19=6", 2) | Out-Null
$d.Content.Find.Execute("12+32=44", $true, $false, $false, $false, $false, $true, 1, $false, "96-57=39", 2) | Out-Null
$d.Content.Find.Execute("25+7=32", $true, $false, $false, $false, $false, $true, 1, $false, "49+19=68", 2) | Out-Null
$d.Content.Find.Execute("63-53=10", $true, $false, $false, $false, $false, $true, 1, $false, "0+96=96", 2) | Out-Null
$d.Content.Find.Execute("93-76=17", $true, $false, $false, $false, $false, $true, 1, $false, "99-86=13", 2) | Out-Null
$d.Content.Find.Execute("18+14=32", $true, $false, $false, $false, $false, $true, 1, $false, "32-30=2", 2) | Out-Null
$d.Content.Find.Execute("9+6=15", $true, $false, $false, $false, $false, $true, 1, $false, "61+20=81", 2) | Out-Null
$d.Content.Find.Execute("52+43=95", $true, $false, $false, $false, $false, $true, 1, $false, "40-5=35", 2) | Out-Null
$d.Content.Find.Execute("12+41=53", $true, $false, $false, $false, $false, $true, 1, $false, "63+3=66", 2) | Out-Null
$d.Content.Find.Execute("45+20=65", $true, $false, $false, $false, $false, $true, 1, $false, "9+88=97", 2) | Out-Null
$d.Content.Find.Execute("86+11=97", $true, $false, $false, $false, $false, $true, 1, $false, "8+2=10", 2) | Out-Null
$d.Content.Find.Execute("4+31=35", $true, $false, $false, $false, $false, $true, 1, $false, "2+77=79", 2) | Out-Null
$d.Content.Find.Execute("59-20=39", $true, $false, $false, $false, $false, $true, 1, $false, "78+17=95", 2) | Out-Null
$d.Content.Find.Execute("40-18=22", $true, $false, $false, $false, $false, $true, 1, $false, "37-28=9", 2) | Out-Null
$d.Content.Find.Execute("63-19=44", $true, $false, $false, $false, $false, $true, 1, $false, "46+28=74", 2) | Out-Null
$d.Content.Find.Execute("58+39=97", $true, $false, $false, $false, $false, $true, 1, $false, "65-41=24", 2) | Out-Null
$d.Content.Find.Execute("65+18=83", $true, $false, $false, $false, $false, $true, 1, $false, "52-6=46", 2) | Out-Null
$d.Content.Find.Execute("22+58=80", $true, $false, $false, $false, $false, $true, 1, $false, "50-30=20", 2) | Out-Null
$d.Content.Find.Execute("73-6=67", $true, $false, $false, $false, $false, $true, 1, $false, "20+63=83", 2) | Out-Null
$d.Content.Find.Execute("54+3=57", $true, $false, $false, $false, $false, $true, 1, $false, "27+6=33", 2) | Out-Null
$d.Content.Find.Execute("55+39=94", $true, $false, $false, $false, $false, $true, 1, $false, "13+24=37", 2) | Out-Null
$d.Content.Find.Execute("69-16=53", $true, $false, $false, $false, $false, $true, 1, $false, "53-38=15", 2) | Out-Null
$d.Content.Find.Execute("13+43=56", $true, $false, $false, $false, $false, $true, 1, $false, "63+16=79", 2) | Out-Null
$d.Content.Find.Execute("87-53=34", $true, $false, $false, $false, $false, $true, 1, $false, "6+46=52", 2) | Out-Null
$d.Content.Find.Execute("70-37=33", $true, $false, $false, $false, $false, $true, 1, $false, "19+7=26", 2) | Out-Null
$d.Content.Find.Execute("43+47=90", $true, $false, $false, $false, $false, $true, 1, $false, "85-13=72", 2) | Out-Null
$d.Content.Find.Execute("50+28=78", $true, $false, $false, $false, $false, $true, 1, $false, "43+19=62", 2) | Out-Null
$d.Content.Find.Execute("59+26=85", $true, $false, $false, $false, $false, $true, 1, $false, "9+88=97", 2) | Out-Null
$d.Content.Find.Execute("42+37=79", $true, $false, $false, $false, $false, $true, 1, $false, "51-49=2", 2) | Out-Null
$d.Content.Find.Execute("84-37=47", $true, $false, $false, $false, $false, $true, 1, $false, "19-11=8", 2) | Out-Null
$d.Content.Find.Execute("49-11=38", $true, $false, $false, $false, $false, $true, 1, $false, "42+43=85", 2) | Out-Null
$d.Content.Find.Execute("71+7=78", $true, $false, $false, $false, $false, $true, 1, $false, "75-60=15", 2) | Out-Null
$d.Content.Find.Execute("7+88=95", $true, $false, $false, $false, $false, $true, 1, $false, "57-55=2", 2) | Out-Null
$d.Content.Find.Execute("26+51=77", $true, $false, $false, $false, $false, $true, 1, $false, "74-45=29", 2) | Out-Null
$d.Content.Find.Execute("61+30=91", $true, $false, $false, $false, $false, $true, 1, $false, "47+49=96", 2) | Out-Null
$d.Content.Find.Execute("75-43=32", $true, $false, $false, $false, $false, $true, 1, $false, "92-55=37", 2) | Out-Null
$d.Content.Find.Execute("18-8=10", $true, $false, $false, $false, $false, $true, 1, $false, "77-11=66", 2) | Out-Null
$d.Content.Find.Execute("51+46=97", $true, $false, $false, $false, $false, $true, 1, $false, "84-16=68", 2) | Out-Null
$d.Content.Find.Execute("53-26=27", $true, $false, $false, $false, $false, $true, 1, $false, "18+63=81", 2) | Out-Null
$d.Content.Find.Execute("83-72=11", $true, $false, $false, $false, $false, $true, 1, $false, "50+13=63", 2) | Out-Null
$d.Content.Find.Execute("96-80=16", $true, $false, $false, $false, $false, $true, 1, $false, "18-16=2", 2) | Out-Null
$d.Content.Find.Execute("31+10=41", $true, $false, $false, $false, $false, $true, 1, $false, "87+4=91", 2) | Out-Null
$d.Content.Find.Execute("0+9=9", $true, $false, $false, $false, $false, $true, 1, $false, "95-0=95", 2) | Out-Null
$d.Content.Find.Execute("83-54=29", $true, $false, $false, $false, $false, $true, 1, $false, "44-24=20", 2) | Out-Null
$d.Content.Find.Execute("2+30=32", $true, $false, $false, $false, $false, $true, 1, $false, "55-28=27", 2) | Out-Null
$d.Content.Find.Execute("90-29=61", $true, $false, $false, $false, $false, $true, 1, $false, "46-22=24", 2) | Out-Null
$d.Content.Find.Execute("16+1=17", $true, $false, $false, $false, $false, $true, 1, $false, "18+41=59", 2) | Out-Null
$d.Content.Find.Execute("44-3=41", $true, $false, $false, $false, $false, $true, 1, $false, "26+49=75", 2) | Out-Null
$d.Content.Find.Execute("90-49=41", $true, $false, $false, $false, $false, $true, 1, $false, "78+14=92", 2) | Out-Null
$d.Content.Find.Execute("34+24=58", $true, $false, $false, $false, $false, $true, 1, $false, "0+20=20", 2) | Out-Null
$d.Content.Find.Execute("88-37=51", $true, $false, $false, $false, $false, $true, 1, $false, "51+9=60", 2) | Out-Null
$d.Content.Find.Execute("13-8=5", $true, $false, $false, $false, $false, $true, 1, $false, "5+68=73", 2) | Out-Null
$d.Content.Find.Execute("85+3=88", $true, $false, $false, $false, $false, $true, 1, $false, "65-38=27", 2) | Out-Null
$d.Content.Find.Execute("21+30=51", $true, $false, $false, $false, $false, $true, 1, $false, "82-25=57", 2) | Out-Null
$d.Content.Find.Execute("28+47=75", $true, $false, $false, $false, $false, $true, 1, $false, "43+33=76", 2) | Out-Null
$d.Content.Find.Execute("93-51=42", $true, $false, $false, $false, $false, $true, 1, $false, "68-47=21", 2) | Out-Null
$d.Content.Find.Execute("26+67=93", $true, $false, $false, $false, $false, $true, 1, $false, "98-17=81", 2) | Out-Null
$d.Content.Find.Execute("23+63=86", $true, $false, $false, $false, $false, $true, 1, $false, "45-28=17", 2) | Out-Null
$d.Content.Find.Execute("28-20=8", $true, $false, $false, $false, $false, $true, 1, $false, "5+10=15", 2) | Out-Null
$d.Content.Find.Execute("17+1=18", $true, $false, $false, $false, $false, $true, 1, $false, "28+51=79", 2) | Out-Null
$d.Content.Find.Execute("27+37=64", $true, $false, $false, $false, $false, $true, 1, $false, "23+26=49", 2) | Out-Null
$d.Content.Find.Execute("65+22=87", $true, $false, $false, $false, $false, $true, 1, $false, "26-15=11", 2) | Out-Null
$d.Content.Find.Execute("76-40=36", $true, $false, $false, $false, $false, $true, 1, $false, "26+35=61", 2) | Out-Null
$d.Content.Find.Execute("51+28=79", $true, $false, $false, $false, $false, $true, 1, $false, "44+28=72", 2) | Out-Null
$d.Content.Find.Execute("62+20=82", $true, $false, $false, $false, $false, $true, 1, $false, "45-5=40", 2) | Out-Null
$d.Content.Find.Execute("8+24=32", $true, $false, $false, $false, $false, $true, 1, $false, "88-21=67", 2) | Out-Null
$d.Content.Find.Execute("9+65=74", $true, $false, $false, $false, $false, $true, 1, $false, "75-46=29", 2) | Out-Null
$d.Content.Find.Execute("79+9=88", $true, $false, $false, $false, $false, $true, 1, $false, "6+84=90", 2) | Out-Null
$d.Content.Find.Execute("53+28=81", $true, $false, $false, $false, $false, $true, 1, $false, "64-26=38", 2) | Out-Null
